$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.8571612750185557, 0.5412043507944918, 1.4583574401596122, 2.6997192601375595)"
$ws.Range("C2").Value = "NIG(1.281384389105698, 0.9927300236049552, 3.38832357365171, 5.073188678935194)"
$ws.Range("D2").Value = "NCT(2.1580483291731913, 1.4355159246686284, -0.5372498754761921, 2.1814039953655016)"
$ws.Range("E2").Value = "JSU(-1.3209888415451365, 1.21751637320558, 2.5746443171940823, 4.80394103860375)"
